$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The merged group A5:A9 ("Type not specified") grows to A6:A11, and a new
# standalone row ("Not Labelled") is inserted above it at row 4.
# Unmerge first so every cell can be written/addressed individually.
$ws.Range("A5:A9").UnMerge()

# Row 2 unchanged: Mixed / Other (Specify) / 3 / 0.6
# Row 3 unchanged: (blank, merged) / Single Cable / 2 / 0.4

# Row 4 (new): Not Labelled / Type not specified / 4 / 0.8
$ws.Range("A4").Value = "Not Labelled"
$ws.Range("B4").Value = "Type not specified"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 0.8

# Row 5 (was row 4): Thermoset / Single Cable / 2 / 0.4
$ws.Range("A5").Value = "Thermoset"
$ws.Range("B5").Value = "Single Cable"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 0.4

# Row 6 (was row 5, now top of the merged group): Type not specified / Single Cable / 9 / 1.7
$ws.Range("A6").Value = "Type not specified"
$ws.Range("B6").Value = "Single Cable"
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 1.7

# Row 7 (was row 6): Multiple Cables Not In Tray or Bundle / 5 / 1
$ws.Range("A7").Value = ""
$ws.Range("B7").Value = "Multiple Cables Not In Tray or Bundle"
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 1

# Row 8 (was row 7): Unknown / 4 / 0.8
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = "Unknown"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 0.8

# Row 9 (was row 8): Other (Specify) / 3 / 0.6
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = "Other (Specify)"
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 0.6

# Row 10 (was row 9): Single Cable Tray / 3 / 0.6
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = "Single Cable Tray"
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 0.6

# Row 11 (new, appended to the bottom of the merged group): Mixed / 1 / 0.2
$ws.Range("A11").Value = ""
$ws.Range("B11").Value = "Mixed"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0.2

# --- Re-merge the "Type not specified" group, now spanning rows 6-11 ---
# Clear formatting first so Merge() does not synthesize partial-border
# style variants for the interior cells of the merged block.
$ws.Range("A6:A11").ClearFormats()
$ws.Range("A6:A11").Merge()

# --- (Re)apply the formatting (bold / bordered / centered-top) used by every
#     other cell in columns A & B. B2 is used as the style donor because it
#     carries the same formatting but was never part of a merged range, so
#     copying from it will not introduce extra partial-border styles. ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:A11").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
